# The commit swaps the raw contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml:
#   before: theme1.xml = "Office Theme" colours, theme2.xml = "Integral" colours
#   after : theme1.xml = "Integral" colours,      theme2.xml = "Office Theme" colours
#
# ppt/theme/theme2.xml is the theme used by the slide master (and by the
# presentation's root theme relationship), so it is what actually drives the
# on-screen look of every slide. We reach it through the modern 12-slot
# DrawingML theme colour scheme (ThemeColorScheme), which maps 1:1 onto
# <a:clrScheme> in document order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink. Driving it through Slide.ThemeColorScheme.Colors(i).RGB edits
# theme2.xml in place without touching anything else (no shape, layout or
# relationship changes), matching the diff's minimal footprint.
#
# RGB() style values below are the little-endian 0xBBGGRR integers that the
# PowerPoint object model's ColorFormat/.RGB setter expects for each of the
# "Office Theme" hex colours that theme2.xml is switched to.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
